$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Manufacturer1075_AT"
$ws.Range("A4").Value = "AuthorisedRep1075_AT"
$ws.Range("A6").Value = "Manufacturer1075_NU"
$ws.Range("A7").Value = "AuthorisedRep1075_NU"
